$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (D text value, E text value). D is $null where the diff does not
# change the Price column for that row. A leading apostrophe forces the
# Price values to stay text (matching the source data, which stores
# numeric-looking prices like "521.76" as text, not as numbers).
$updates = @(
  @{ Row = 2;  D = "60.632.78";  E = "  +3.38%  " },
  @{ Row = 3;  D = "2.678.29";   E = "  +2.12%  " },
  @{ Row = 4;  D = $null;        E = "  -0.14%  " },
  @{ Row = 5;  D = "521.76";     E = "  +2.22%  " },
  @{ Row = 6;  D = "146.98";     E = "  +2.32%  " },
  @{ Row = 7;  D = $null;        E = "  +0.02%  " },
  @{ Row = 8;  D = "0.579";      E = "  +2.53%  " },
  @{ Row = 9;  D = "2.703.51";   E = "  +1.82%  " },
  @{ Row = 10; D = "6.46";       E = "  +1.45%  " },
  @{ Row = 11; D = $null;        E = "  +1.19%  " },
  @{ Row = 12; D = "0.341";      E = "  +1.47%  " },
  @{ Row = 13; D = $null;        E = "  +1.48%  " },
  @{ Row = 14; D = "3.154.29";   E = "  +1.96%  " },
  @{ Row = 15; D = "60.621.94";  E = "  +3.21%  " },
  @{ Row = 16; D = "21.39";      E = "  +1.91%  " },
  @{ Row = 17; D = "2.771.79";   E = "  +4.53%  " },
  @{ Row = 18; D = $null;        E = "  +1.80%  " },
  @{ Row = 19; D = "352.66";     E = "  +2.59%  " },
  @{ Row = 20; D = "4.57";       E = "  +0.60%  " },
  @{ Row = 21; D = "10.55";      E = "  +2.01%  " },
  @{ Row = 22; D = "6.35";       E = "  +4.38%  " },
  @{ Row = 23; D = "0.997";      E = "  -0.04%  " },
  @{ Row = 24; D = "62.89";      E = "  +2.90%  " },
  @{ Row = 25; D = "0.425";      E = "  +1.10%  " },
  @{ Row = 26; D = "0.169";      E = "  +4.63%  " },
  @{ Row = 27; D = "0.995";      E = "  -0.18%  " },
  @{ Row = 28; D = "0.0₃0820";   E = "  +2.14%  " },
  @{ Row = 29; D = "7.29";       E = "  +2.79%  " },
  @{ Row = 30; D = "6.90";       E = "  +7.20%  " },
  @{ Row = 31; D = $null;        E = "  +0.09%  " },
  @{ Row = 32; D = $null;        E = "  +1.70%  " },
  @{ Row = 33; D = "19.12";      E = "  +1.36%  " },
  @{ Row = 34; D = "149.27";     E = "  -0.60%  " },
  @{ Row = 35; D = "4.33";       E = "  +7.85%  " },
  @{ Row = 36; D = "0.949";      E = "  -7.71%  " },
  @{ Row = 37; D = "1.23";       E = "  +6.70%  " },
  @{ Row = 38; D = $null;        E = "  +11.39%  " },
  @{ Row = 39; D = "0.876";      E = "  +2.91%  " },
  @{ Row = 40; D = "36.69";      E = "  +0.65%  " },
  @{ Row = 41; D = "3.71";       E = "  +0.74%  " },
  @{ Row = 42; D = "284.69";     E = "  +2.08%  " },
  @{ Row = 43; D = "20.08";      E = "  +3.04%  " },
  @{ Row = 44; D = "0.0993";     E = "  +0.98%  " },
  @{ Row = 45; D = $null;        E = "  -0.24%  " },
  @{ Row = 46; D = $null;        E = "  -0.14%  " },
  @{ Row = 47; D = "2.145.67";   E = "  +8.66%  " },
  @{ Row = 48; D = "0.0542";     E = "  +1.11%  " },
  @{ Row = 49; D = "4.87";       E = "  +4.25%  " },
  @{ Row = 50; D = "0.0235";     E = "  +3.09%  " },
  @{ Row = 51; D = "10.46";      E = "  +1.82%  " }
)

foreach ($u in $updates) {
  $r = $u.Row
  if ($null -ne $u.D) {
    $ws.Cells.Item($r, 4).Value = "'" + $u.D
  }
  $ws.Cells.Item($r, 5).Value = $u.E
}
